$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 452 - this shifts the existing rows
# 452..493 down to 453..494, preserving all of their original data (which is
# exactly what the target diff shows row-by-row).
$ws.Rows.Item(452).Insert()

# Populate the newly inserted row 452 with the new price-report record.
$ws.Cells.Item(452, 1).Value = 3
$ws.Cells.Item(452, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(452, 3).Value = "Coquimbo"
$ws.Cells.Item(452, 4).Value = 45132
$ws.Cells.Item(452, 5).Value = 5
$ws.Cells.Item(452, 6).Value = 100112001
$ws.Cells.Item(452, 7).Value = "Berenjena"
$ws.Cells.Item(452, 8).Value = "Sin especificar"
$ws.Cells.Item(452, 9).Value = "Primera"
$ws.Cells.Item(452, 10).Value = 90
$ws.Cells.Item(452, 11).Value = 7000
$ws.Cells.Item(452, 12).Value = 7500
$ws.Cells.Item(452, 13).Value = 7222
$ws.Cells.Item(452, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(452, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(452, 16).Value = 120
$ws.Cells.Item(452, 17).Value = 60
$ws.Cells.Item(452, 18).Value = "Hortaliza"

# Match the date-number style used by column D in the rest of the sheet.
$ws.Cells.Item(452, 4).NumberFormat = $ws.Cells.Item(453, 4).NumberFormat
